# "Guardar archivos a Excel"
# The second row of the SearchProduct sheet (search result #2) is replaced:
# it used to describe a Televisor (TV) search result and now describes a
# "Balon" (soccer ball) search result instead. The workbook selection is
# also moved from A2 to C9.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SearchProduct")

# Row 2: NumProduct / product search term / validation/detail text
$ws.Range("A2").Value = "Balon"
$ws.Range("B2").Value = "2"
$ws.Range("C2").Value = "GOLTY`nBALÓN DE FÚTBOL PARA NIÑOS GOLTY DINO No4...`nPor E & M S.a"

# Keep the row's height on "auto" (matches how the sheet looked before the
# edit) instead of leaving it pinned to the height Excel computes right
# after assigning a multi-line value.
$ws.Rows(2).AutoFit()

# The user had ended up with C9 selected when the file was saved.
$ws.Range("C9").Select()

$wb.Save()
